$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "56.807.44"
$ws.Range("E2").Value = "  +2.24%  "

# Row 3
$ws.Range("D3").Value = "3.001.51"

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").Value = "'512.49"
$ws.Range("E5").Value = "  +4.72%  "

# Row 6
$ws.Range("D6").Value = "'138.54"

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").Value = "'0.435"
$ws.Range("E8").Value = "  +3.91%  "

# Row 9
$ws.Range("D9").Value = "'7.49"
$ws.Range("E9").Value = "  +4.45%  "

# Row 10
$ws.Range("E10").Value = "  +7.22%  "

# Row 11
$ws.Range("D11").Value = "'0.357"
$ws.Range("E11").Value = "  +3.24%  "

# Row 12
$ws.Range("E12").Value = "  +1.80%  "

# Row 13
$ws.Range("D13").Value = "3.525.49"
$ws.Range("E13").Value = "  +1.60%  "

# Row 14
$ws.Range("D14").Value = "'25.82"
$ws.Range("E14").Value = "  +4.64%  "

# Row 15
$ws.Range("D15").Value = "'0.0000156"
$ws.Range("E15").Value = "  +11.72%  "

# Row 16
$ws.Range("D16").Value = "56.939.37"
$ws.Range("E16").Value = "  +2.52%  "

# Row 17
$ws.Range("D17").Value = "3.007.20"
$ws.Range("E17").Value = "  +1.10%  "

# Row 18
$ws.Range("D18").Value = "'5.92"
$ws.Range("E18").Value = "  +5.24%  "

# Row 19
$ws.Range("D19").Value = "'12.52"
$ws.Range("E19").Value = "  +2.83%  "

# Row 20
$ws.Range("D20").Value = "'7.83"
$ws.Range("E20").Value = "  +5.13%  "

# Row 21
$ws.Range("D21").Value = "'326.29"
$ws.Range("E21").Value = "  +3.03%  "

# Row 22
$ws.Range("E22").Value = "  +0.02%  "

# Row 23
$ws.Range("D23").Value = "'0.486"
$ws.Range("E23").Value = "  +5.59%  "

# Row 24
$ws.Range("D24").Value = "'63.48"
$ws.Range("E24").Value = "  +5.55%  "

# Row 25
$ws.Range("D25").Value = "'0.171"
$ws.Range("E25").Value = "  +5.85%  "

# Row 26
$ws.Range("E26").Value = "  -0.50%  "

# Row 27
$ws.Range("D27").Value = "0.0₃0912"
$ws.Range("E27").Value = "  +9.10%  "

# Row 28
$ws.Range("D28").Value = "'6.64"
$ws.Range("E28").Value = "  +2.81%  "

# Row 29
$ws.Range("D29").Value = "'7.04"
$ws.Range("E29").Value = "  +8.23%  "

# Row 30
$ws.Range("E30").Value = "  +6.09%  "

# Row 31
$ws.Range("D31").Value = "'1.81"
$ws.Range("E31").Value = "  +6.84%  "

# Row 32
$ws.Range("D32").Value = "'20.54"
$ws.Range("E32").Value = "  +5.80%  "

# Row 33
$ws.Range("D33").Value = "'156.88"
$ws.Range("E33").Value = "  +4.82%  "

# Row 34
$ws.Range("D34").Value = "'4.56"
$ws.Range("E34").Value = "  +4.80%  "

# Row 35
$ws.Range("D35").Value = "'5.69"
$ws.Range("E35").Value = "  +0.62%  "

# Row 36
$ws.Range("D36").Value = "'1.27"
$ws.Range("E36").Value = "  -1.88%  "

# Row 37
$ws.Range("D37").Value = "'0.0678"
$ws.Range("E37").Value = "  +4.82%  "

# Row 38
$ws.Range("D38").Value = "'23.86"
$ws.Range("E38").Value = "  +3.48%  "

# Row 39
$ws.Range("D39").Value = "3.042.15"
$ws.Range("E39").Value = "  +1.53%  "

# Row 40
$ws.Range("D40").Value = "'36.98"
$ws.Range("E40").Value = "  +1.37%  "

# Row 41
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.34%  "

# Row 42
$ws.Range("D42").Value = "2.292.08"
$ws.Range("E42").Value = "  +7.82%  "

# Row 43
$ws.Range("D43").Value = "'0.648"

# Row 44
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.05%  "

# Row 45
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "'1.42"
$ws.Range("E45").Value = "  +3.38%  "

# Row 46
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").Value = "'3.67"
$ws.Range("E46").Value = "  +4.44%  "

# Row 47
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0240"
$ws.Range("E47").Value = "  +2.85%  "

# Row 48
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").Value = "'1.94"
$ws.Range("E48").Value = "  +9.42%  "

# Row 49
$ws.Range("D49").Value = "'5.88"
$ws.Range("E49").Value = "  +6.31%  "

# Row 50
$ws.Range("D50").Value = "'19.31"
$ws.Range("E50").Value = "  +1.24%  "

# Row 51
$ws.Range("E51").Value = "  +5.26%  "
